$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Intro_0")
$ws2 = $wb.Worksheets.Item("BD_1")

# Clear the postcode value in BD_1!I2 (was "E13 6SE") - "Customised Timestamp code" edit
[void]$ws2.Range("I2").ClearContents()

# Update the selection remembered on the Intro_0 sheet (it is no longer the active tab)
[void]$ws1.Range("I5").Select()

# Make BD_1 the active sheet and leave the selection on I2 (the cleared cell)
[void]$ws2.Activate()
[void]$ws2.Range("I2").Select()
